$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:I1), mirroring the style of the existing D1/E1 headers
$ws.Range("F1").Value = "09-04-2025 Status"
$ws.Range("G1").Value = "09-04-2025 Time"
$ws.Range("H1").Value = "10-04-2025 Status"
$ws.Range("I1").Value = "10-04-2025 Time"

# Copy the header style from D1 (existing bold/centered/bordered header style) to the new headers
$ws.Range("D1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122) | Out-Null

# Restore the values (PasteSpecial formats may have overwritten text with source D1 text)
$ws.Range("F1").Value = "09-04-2025 Status"
$ws.Range("G1").Value = "09-04-2025 Time"
$ws.Range("H1").Value = "10-04-2025 Status"
$ws.Range("I1").Value = "10-04-2025 Time"

# Fill data rows 2-23 for the new columns, matching the existing D/E pattern (A / 00:00:00)
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 6).Value = "A"
    $ws.Cells.Item($r, 7).Value = "00:00:00"
    $ws.Cells.Item($r, 8).Value = "A"
    $ws.Cells.Item($r, 9).Value = "00:00:00"
}
